{"js": "const pairs = [\n  [\"2023-11-12 Sunday\", \"2023-11-13 Monday\"],\n  [\"29\u00d781=2349\", \"51\u00d773=3723\"],\n  [\"90\u00d770=6300\", \"66\u00d765=4290\"],\n  [\"69\u00d748=3312\", \"93\u00d781=7533\"],\n  [\"83\u00d787=7221\", \"30\u00d776=2280\"],\n  [\"23\u00d796=2208\", \"87\u00d713=1131\"],\n  [\"32\u00d794=3008\", \"39\u00d778=3042\"],\n  [\"79\u00d779=6241\", \"20\u00d779=1580\"],\n  [\"92\u00d720=1840\", \"59\u00d775=4425\"],\n  [\"96\u00d767=6432\", \"91\u00d783=7553\"],\n  [\"50\u00d785=4250\", \"17\u00d724=408\"],\n  [\"46\u00d763=2898\", \"99\u00d749=4851\"],\n  [\"90\u00d718=1620\", \"12\u00d783=996\"],\n  [\"67\u00d723=1541\", \"21\u00d732=672\"],\n  [\"35\u00d754=1890\", \"23\u00d772=1656\"],\n  [\"95\u00d720=1900\", \"92\u00d729=2668\"],\n  [\"11\u00d734=374\", \"25\u00d767=1675\"],\n  [\"95\u00d745=4275\", \"68\u00d780=5440\"],\n  [\"72\u00d723=1656\", \"25\u00d711=275\"],\n  [\"76\u00d777=5852\", \"69\u00d714=966\"],\n  [\"33\u00d743=1419\", \"11\u00d723=253\"],\n  [\"62\u00d747=2914\", \"45\u00d715=675\"],\n  [\"89\u00d781=7209\", \"84\u00d739=3276\"],\n  [\"69\u00d752=3588\", \"94\u00d788=8272\"],\n  [\"46\u00d796=4416\", \"76\u00d718=1368\"],\n  [\"42\u00d711=462\", \"66\u00d785=5610\"],\n];\n\nconst body = context.document.body;\nfor (const [oldText, newText] of pairs) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n  for (const item of results.items) {\n    item.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n\n$pairs = @(\n    @('2023-11-12 Sunday', '2023-11-13 Monday'),\n    @('29\u00d781=2349', '51\u00d773=3723'),\n    @('90\u00d770=6300', '66\u00d765=4290'),\n    @('69\u00d748=3312', '93\u00d781=7533'),\n    @('83\u00d787=7221', '30\u00d776=2280'),\n    @('23\u00d796=2208', '87\u00d713=1131'),\n    @('32\u00d794=3008', '39\u00d778=3042'),\n    @('79\u00d779=6241', '20\u00d779=1580'),\n    @('92\u00d720=1840', '59\u00d775=4425'),\n    @('96\u00d767=6432', '91\u00d783=7553'),\n    @('50\u00d785=4250', '17\u00d724=408'),\n    @('46\u00d763=2898', '99\u00d749=4851'),\n    @('90\u00d718=1620', '12\u00d783=996'),\n    @('67\u00d723=1541', '21\u00d732=672'),\n    @('35\u00d754=1890', '23\u00d772=1656'),\n    @('95\u00d720=1900', '92\u00d729=2668'),\n    @('11\u00d734=374', '25\u00d767=1675'),\n    @('95\u00d745=4275', '68\u00d780=5440'),\n    @('72\u00d723=1656', '25\u00d711=275'),\n    @('76\u00d777=5852', '69\u00d714=966'),\n    @('33\u00d743=1419', '11\u00d723=253'),\n    @('62\u00d747=2914', '45\u00d715=675'),\n    @('89\u00d781=7209', '84\u00d739=3276'),\n    @('69\u00d752=3588', '94\u00d788=8272'),\n    @('46\u00d796=4416', '76\u00d718=1368'),\n    @('42\u00d711=462', '66\u00d785=5610'),\n)\n\nforeach ($pair in $pairs) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $pair[0]\n    $find.Replacement.Text = $pair[1]\n    $find.Execute([ref]$pair[0], $false, $true, $false, $false, $false, $true, 1, $false, [ref]$pair[1], 2)\n}\n"}
